# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") previously held the raw strike-out
# count (Strike#). This regenerates that column with the corrected
# per-game strikeout ("K") values for rows 2-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 35 (in row order).
$kValues = @(1, 4, 5, 2, 5, 0, 4, 4, 3, 8, 5, 5, 10, 2, 3, 7, 5, 1, 6, 6, 3, 5, 4, 5, 2, 8, 2, 7, 6, 3, 4, 5, 3, 4)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
